$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(203957296, Omri Ben Shabat: -8,7)"
$ws.Range("B1").Value = "(206532695, Matan Vakrat: 8,-6)"
$ws.Range("C1").Value = "(302962915, Asher  Odeh: -8,-6)"
$ws.Range("D1").Value = "(308035542, Anastasia  Kubi: -8,-1)"
$ws.Range("E1").Value = "(311177802, Christina  Uksusman: 5,6)"
$ws.Range("F1").Value = "(305251175, Or  Leder: -9,-4)"

$ws.Range("A3").Value = "cost: 834.1838832406117"
$ws.Range("A4").Value = "time: 101.14798540507647"
